# Apply automatic update to rows 2, 3 and 5 of the active worksheet.
# The update rotates the species/observation data among rows 2, 3 and 5
# (row2 <- old row3 data, row3 <- old row5 data, row5 <- old row2 data),
# while the "Taxonsorteringsordning" (column B) values are bumped by 14
# relative to the donor row's original value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets updated with the (shifted) data previously in row 3
$ws.Range("A2").Value = 112000857
$ws.Range("B2").Value = 77650
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 766907
$ws.Range("R2").Value = 7097936
$ws.Range("Z2").Value = "10:12"
$ws.Range("AB2").Value = "10:12"

# Row 3 gets updated with the (shifted) data previously in row 5
$ws.Range("A3").Value = 112003034
$ws.Range("B3").Value = 90814
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 4364
$ws.Range("F3").Value = "Dropptaggsvamp"
$ws.Range("G3").Value = "Hydnellum ferrugineum"
$ws.Range("H3").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q3").Value = 766713
$ws.Range("R3").Value = 7098084
$ws.Range("S3").Value = 25

# Row 5 gets updated with the (shifted) data previously in row 2
$ws.Range("A5").Value = 112000795
$ws.Range("B5").Value = 89571
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 5432
$ws.Range("F5").Value = "Granticka"
$ws.Range("G5").Value = "Porodaedalea chrysoloma"
$ws.Range("H5").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q5").Value = 766898
$ws.Range("R5").Value = 7097941
$ws.Range("S5").Value = 100
$ws.Range("Z5").Value = "10:09"
$ws.Range("AB5").Value = "10:09"
